$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.147.09'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '2.210.96'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'228.83"
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").Value = "'0.630"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = "'64.12"
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = "'0.398"
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").Value = "'0.0862"
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("E11").Value = '  +0.66%  '
$ws.Range("D12").Value = "'16.16"
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").Value = '2.537.32'
$ws.Range("E13").Value = '  +2.15%  '
$ws.Range("D14").Value = "'22.27"
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("D16").Value = "'5.62"
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '2.203.59'
$ws.Range("E17").Value = '  +1.44%  '
$ws.Range("D18").Value = '40.063.23'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("D19").Value = '0.0₃0913'
$ws.Range("E19").Value = '  +6.85%  '
$ws.Range("D20").Value = "'72.40"
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").Value = "'6.12"
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").Value = "'232.21"
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").Value = "'9.70"
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").Value = "'171.91"
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("D28").Value = "'0.141"
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").Value = "'20.18"
$ws.Range("E30").Value = '  +2.54%  '
$ws.Range("D31").Value = "'2.74"
$ws.Range("E31").Value = '  +5.85%  '
$ws.Range("E32").Value = '  +1.16%  '
$ws.Range("D33").Value = "'4.61"
$ws.Range("E33").Value = '  -1.54%  '
$ws.Range("D34").Value = "'4.75"
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("D35").Value = "'7.07"
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = "'0.0625"
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("D37").Value = "'3.89"
$ws.Range("E37").Value = '  +9.15%  '
$ws.Range("D38").Value = "'2.48"
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("B39").Value = 'BinanceUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("B40").Value = 'FTXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D40").Value = "'4.98"
$ws.Range("E40").Value = '  +17.65%  '
$ws.Range("D41").Value = "'104.05"
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").Value = "'17.89"
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("E44").Value = '  +3.46%  '
$ws.Range("D45").Value = '1.524.67'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = "'8.26"
$ws.Range("E46").Value = '  +6.11%  '
$ws.Range("E47").Value = '  +0.75%  '
$ws.Range("D48").Value = "'0.0929"
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").Value = "'0.000196"
$ws.Range("E50").Value = '  +34.08%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.416.20'
$ws.Range("E51").Value = '  +2.07%  '